# Generate Report for Handback
#
# For the "27384982-d757-4886-a685-700753facc32" file row (row 6) on both the
# "zh-cn" and "de-de" sheets, the handback report now has a generated
# Latest Target File / Latest Handback File / Latest Handback DateTime, plus
# an Error Detail describing that the handback file version is stale. The
# Error Detail column (P) is widened to fit the long message.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b1993d9f41093dce23d03efb122f1f1ca5add56b/e2e/27384982-d757-4886-a685-700753facc32.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1468bd9a1fe9d1f3678717b693da69d0a076e259/e2e/27384982-d757-4886-a685-700753facc32.md."
$latestUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1468bd9a1fe9d1f3678717b693da69d0a076e259/e2e/27384982-d757-4886-a685-700753facc32.md"
$handbackName = "27384982-d757-4886-a685-700753facc32.md"

function Update-HandbackRow($sheetName, $handbackDateTime) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Latest Target File (J6) mirrors the Latest Handoff File (G6) once a
    # handback round-trip has happened.
    $targetFile = $ws.Cells.Item(6, 7).Value2
    $ws.Cells.Item(6, 10).Value = $targetFile

    # Latest Handback File (I6): the handback markdown file name, rendered as
    # a hyperlink just like column A's entry.
    $i6 = $ws.Cells.Item(6, 9)
    $i6.Value = $handbackName
    $ws.Hyperlinks.Add($i6, $latestUrl, [Type]::Missing, [Type]::Missing, $handbackName) | Out-Null
    # Restore the workbook's existing custom Hyperlink look (underlined,
    # explicit blue) instead of the theme-based one `Hyperlinks.Add` applies.
    $i6.Font.Name = "Calibri"
    $i6.Font.Underline = $true
    $i6.Font.Color = 15570276

    # Latest Handback DateTime (K6).
    $ws.Cells.Item(6, 11).Value = $handbackDateTime

    # Error Detail (P6).
    $ws.Cells.Item(6, 16).Value = $errorDetail

    # Widen the Error Detail column so the long message is readable.
    $ws.Columns.Item(16).ColumnWidth = 39.17
}

Update-HandbackRow "zh-cn" "2016-08-16 02:39:30"
Update-HandbackRow "de-de" "2016-08-16 02:39:38"
